# edit.ps1 - apply the "participant_event_list2" commit to the deck:
#   - refresh the cached datetimeFigureOut field (17/04/2019 -> 24/04/2019)
#     everywhere it appears (slide master + every slide layout)
#   - reword the "Ete 2019" / "Hiver 2019" event-card headers
#   - reword the intro/outro instructions paragraphs
#   - replace the "Netflix Friends" placeholder event with a real one
#   - reword the "Paris Express" card title to "Grand Paris"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder field: 17/04/2019 -> 24/04/2019, in the slide master
#    and in every custom (slide) layout.
# ---------------------------------------------------------------------------
function Update-DateField {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "17/04/2019") {
                $tr.Text = "24/04/2019"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-DateField $master.CustomLayouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# Helper: retype a sub-range of a paragraph's text so the saved XML keeps
# (or creates) a dedicated run for that slice, mirroring how PowerPoint
# splits runs as a user retypes part of a line.
# ---------------------------------------------------------------------------
function Set-Chars {
    param($para, [int]$start, [int]$length, [string]$text)
    $rng = $para.Characters($start, $length)
    $rng.Text = $text
}

# ---------------------------------------------------------------------------
# Slide 1
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# "Ete 2019" -> "Fin 2019" (split across 3 runs: "Fin" / " " / "2019")
$rect4 = $s1.Shapes.Item(2)
$para = $rect4.TextFrame2.TextRange.Paragraphs(1, 1)
Set-Chars $para 1 3 "Fin"
Set-Chars $para 5 4 "2019"

# "Hiver 2019" -> "Debut 2019" (split across 3 runs: "D" / "ebut " / "2019")
$rect6 = $s1.Shapes.Item(4)
$para = $rect6.TextFrame2.TextRange.Paragraphs(1, 1)
Set-Chars $para 1 1 "D"
Set-Chars $para 2 5 "ébut "
Set-Chars $para 7 4 "2019"

# Intro / outro paragraphs, reworded
$rect13 = $s1.Shapes.Item(11)
$tr13 = $rect13.TextFrame2.TextRange

$introPara = $tr13.Paragraphs(3, 1)
$introPara.Text = "Vous trouverez ici la liste des événements dont nous vous invitons à prendre connaissance pour l’expérience. Les mots entre guillemets sont ceux qui apparaîtront durant l’expérience."

$outroPara = $tr13.Paragraphs(9, 1)
$outroPara.Text = "Votre connaissance de ces événements et de leurs dates est fondamentale pour le bon déroulement de l’expérience. Bon courage!"

# ---------------------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# "Hiver 2019" / "Netflix Friends" fictitious card -> real "viandes artificielles" card
$rect6b = $s2.Shapes.Item(4)
$tr6b = $rect6b.TextFrame2.TextRange

$tr6b.Paragraphs(1, 1).Text = "2020"
$tr6b.Paragraphs(2, 1).Text = "« viandes artificielles »"
$tr6b.Paragraphs(3, 1).Text = "Réel"

$descPara = $tr6b.Paragraphs(4, 1)
$descFull = "Commercialisation de viandes artificielles"
Set-Chars $descPara 1 ($descPara.Length - 1) $descFull
Set-Chars $descPara ($descPara.Length - 1 - 13) 13 "artificielles"

# "« Paris Express »" -> "« Grand Paris »" (split across 3 runs: "«" / " Grand Paris" / " »")
$rect18 = $s2.Shapes.Item(16)
$para18 = $rect18.TextFrame2.TextRange.Paragraphs(2, 1)
Set-Chars $para18 1 1 "«"
Set-Chars $para18 2 ($para18.Length - 1 - 2) " Grand Paris"
Set-Chars $para18 ($para18.Length - 1 - 1) 2 " »"

Write-Output "edit complete"
